$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / date text updates (shared-string runs collapse to plain text; formatting unaffected) ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Style-class transitions: "N/A" (text, style 14) <-> numeric (style 15) ---
# C15: was "N/A" (shared string), becomes numeric 1 -> borrow numeric format from D15
$ws.Range("C15").Value = 1
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# C26: was "N/A" (shared string), becomes numeric 1 -> borrow numeric format from D26
$ws.Range("C26").Value = 1
$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# C27: was numeric 4, becomes "N/A" text -> force text via leading apostrophe, then
# borrow the "N/A" text format (style 14, General) from C23 afterwards so the
# resulting style/number-format matches the other "N/A" cells exactly.
$ws.Range("C27").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Remaining numeric updates (style class unchanged) ---
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = -13.636363636363
$ws.Range("L15").Value = -17.391304347826
$ws.Range("M15").Value = 46.153846153846
$ws.Range("N15").Value = 35.714285714285
$ws.Range("C16").Value = 14
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 50
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 78.571428571428
$ws.Range("I16").Value = 249
$ws.Range("J16").Value = 222
$ws.Range("K16").Value = 12.162162162162
$ws.Range("L16").Value = 76.595744680851
$ws.Range("M16").Value = 23.267326732673
$ws.Range("N16").Value = -73.92670157068
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 72
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = 38.461538461538
$ws.Range("I17").Value = 439
$ws.Range("J17").Value = 317
$ws.Range("K17").Value = 38.485804416403
$ws.Range("L17").Value = 98.642533936651
$ws.Range("M17").Value = 174.375
$ws.Range("N17").Value = 43.464052287581
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 10
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 102
$ws.Range("K18").Value = 15.686274509803
$ws.Range("L18").Value = 14.563106796116
$ws.Range("M18").Value = -30.994152046783
$ws.Range("N18").Value = -91.486291486291
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 58.823529411764
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -4.109589041095
$ws.Range("I19").Value = 630
$ws.Range("J19").Value = 699
$ws.Range("K19").Value = -9.871244635193
$ws.Range("L19").Value = 131.617647058824
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = -12.133891213389
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 138.461538461538
$ws.Range("I20").Value = 193
$ws.Range("J20").Value = 110
$ws.Range("K20").Value = 75.454545454545
$ws.Range("L20").Value = 98.969072164948
$ws.Range("M20").Value = 96.938775510204
$ws.Range("N20").Value = -85.142417244033
$ws.Range("C21").Value = 76
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 80.95238095238
$ws.Range("F21").Value = 237
$ws.Range("G21").Value = 183
$ws.Range("H21").Value = 29.508196721311
$ws.Range("I21").Value = 1650
$ws.Range("J21").Value = 1475
$ws.Range("K21").Value = 11.864406779661
$ws.Range("L21").Value = 91.860465116279
$ws.Range("M21").Value = 71.696149843912
$ws.Range("N21").Value = -64.871194379391
$ws.Range("C22").Value = 3
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 31
$ws.Range("K22").Value = 72.222222222222
$ws.Range("L22").Value = 138.461538461538
$ws.Range("M22").Value = 63.157894736842
$ws.Range("C24").Value = 54
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = 3.846153846153
$ws.Range("F24").Value = 231
$ws.Range("G24").Value = 190
$ws.Range("H24").Value = 21.578947368421
$ws.Range("I24").Value = 1666
$ws.Range("J24").Value = 1395
$ws.Range("K24").Value = 19.426523297491
$ws.Range("L24").Value = 67.605633802816
$ws.Range("M24").Value = 65.277777777777
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -15.384615384615
$ws.Range("F25").Value = 96
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = 26.315789473684
$ws.Range("I25").Value = 692
$ws.Range("J25").Value = 527
$ws.Range("K25").Value = 31.309297912713
$ws.Range("L25").Value = 55.855855855855
$ws.Range("M25").Value = 80.678851174934
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("I26").Value = 28
$ws.Range("J26").Value = 35
$ws.Range("K26").Value = -20
$ws.Range("L26").Value = -20
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 85.714285714285
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 21.428571428571
$ws.Range("L27").Value = 49.122807017543
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = -14.285714285714
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = -84.210526315789
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 150
$ws.Range("L29").Value = -16.666666666666
$ws.Range("M29").Value = 66.666666666666
$ws.Range("N29").Value = -85.294117647058

